# Updated cryptos list on Sun Jan 14 14:47:12 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    # Force the cell to be treated as text so that numeric-looking strings
    # (e.g. "304.40", "2.537.72") are not auto-converted to numbers, then
    # restore the cell's original (default) style so no stray formatting
    # is introduced.
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "42.903.51"
$ws.Range("E2").Value = "  +0.14%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.537.72"
$ws.Range("E3").Value = "  -0.18%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.07%  "

# Row 5 - BNB
Set-TextValue "D5" "304.40"
$ws.Range("E5").Value = "  +1.86%  "

# Row 6 - Solana
Set-TextValue "D6" "98.90"
$ws.Range("E6").Value = "  +7.66%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  +0.69%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.09%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  -0.20%  "

# Row 10 - Avalanche
Set-TextValue "D10" "37.05"
$ws.Range("E10").Value = "  +2.98%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  +2.50%  "

# Row 12 - Polkadot
$ws.Range("E12").Value = "  +1.16%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  -1.50%  "

# Row 14 - WrappedliquidstakedEther2.0
Set-TextValue "D14" "2.927.26"
$ws.Range("E14").Value = "  -0.09%  "

# Row 15 - WrappedEther
Set-TextValue "D15" "2.585.39"
$ws.Range("E15").Value = "  -0.01%  "

# Row 16 - Chainlink
Set-TextValue "D16" "15.14"
$ws.Range("E16").Value = "  +6.74%  "

# Row 17 - Polygon
Set-TextValue "D17" "0.876"
$ws.Range("E17").Value = "  -0.04%  "

# Row 18 - WrappedBTC
Set-TextValue "D18" "42.930.39"
$ws.Range("E18").Value = "  +0.24%  "

# Row 19 - InternetComputer(DFINITY)
Set-TextValue "D19" "13.11"
$ws.Range("E19").Value = "  +3.82%  "

# Row 20 - ShibaInu
$ws.Range("E20").Value = "  +1.07%  "

# Row 21 - Uniswap
Set-TextValue "D21" "6.56"
$ws.Range("E21").Value = "  +0.45%  "

# Row 22 - Litecoin
Set-TextValue "D22" "71.69"
$ws.Range("E22").Value = "  -0.08%  "

# Row 23 - BitcoinCash
Set-TextValue "D23" "253.76"
$ws.Range("E23").Value = "  -0.49%  "

# Row 24 - PancakeSwap
$ws.Range("E24").Value = "  +0.30%  "

# Row 25 - ImmutableX
Set-TextValue "D25" "2.06"
$ws.Range("E25").Value = "  -3.08%  "

# Row 26 - EthereumClassic
Set-TextValue "D26" "27.73"
$ws.Range("E26").Value = "  -4.07%  "

# Row 27 - Dai
Set-TextValue "D27" "0.999"
$ws.Range("E27").Value = "  -0.18%  "

# Row 28 - Toncoin
Set-TextValue "D28" "2.33"
$ws.Range("E28").Value = "  +10.39%  "

# Row 29 - now InjectiveProtocol (was Cosmos)
$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D29" "39.24"
$ws.Range("E29").Value = "  +7.44%  "

# Row 30 - now Cosmos (was InjectiveProtocol)
$ws.Range("B30").Value = "Cosmos"
$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D30" "10.16"
$ws.Range("E30").Value = "  -0.43%  "

# Row 31 - Filecoin
Set-TextValue "D31" "6.22"
$ws.Range("E31").Value = "  +2.95%  "

# Row 32 - Monero
Set-TextValue "D32" "157.34"
$ws.Range("E32").Value = "  +3.16%  "

# Row 33 - ARBITRUM
$ws.Range("E33").Value = "  +0.14%  "

# Row 34 - Celestia
Set-TextValue "D34" "19.17"
$ws.Range("E34").Value = "  +8.84%  "

# Row 35 - Hedera
$ws.Range("E35").Value = "  +0.93%  "

# Row 36 - LidoDAOToken
$ws.Range("E36").Value = "  -2.15%  "

# Row 37 - WEMIXToken
Set-TextValue "D37" "2.62"
$ws.Range("E37").Value = "  -4.80%  "

# Row 38 - Kaspa
$ws.Range("E38").Value = "  +1.14%  "

# Row 39 - EnergySwap
Set-TextValue "D39" "24.66"
$ws.Range("E39").Value = "  +6.86%  "

# Row 40 - Stellar
$ws.Range("E40").Value = "  +0.75%  "

# Row 41 - ApeXProtocol
$ws.Range("E41").Value = "  +8.69%  "

# Row 42 - NEARProtocol
Set-TextValue "D42" "3.44"
$ws.Range("E42").Value = "  +0.96%  "

# Row 43 - RenderToken
Set-TextValue "D43" "3.91"
$ws.Range("E43").Value = "  +2.05%  "

# Row 44 - now VeChain (was Maker)
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D44" "0.0305"
$ws.Range("E44").Value = "  -1.38%  "

# Row 45 - now Maker (was VeChain)
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D45" "2.082.19"
$ws.Range("E45").Value = "  -0.34%  "

# Row 46 - FirstDigitalUSD
Set-TextValue "D46" "0.998"
$ws.Range("E46").Value = "  +0.02%  "

# Row 47 - BitcoinSV
Set-TextValue "D47" "86.35"
$ws.Range("E47").Value = "  +2.50%  "

# Row 48 - FraxShare
Set-TextValue "D48" "9.00"
$ws.Range("E48").Value = "  -1.23%  "

# Row 49 - RocketPoolETH
Set-TextValue "D49" "2.785.08"
$ws.Range("E49").Value = "  -0.01%  "

# Row 50 - ordi
Set-TextValue "D50" "73.69"
$ws.Range("E50").Value = "  +6.70%  "

# Row 51 - Algorand
$ws.Range("E51").Value = "  +2.04%  "
